# "first run with logo"
# - widen column A
# - flip F column to 1 for a set of rows (marking isMCSB / "has logo")
# - update the sheet scroll/selection to match the author's final cursor spot

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A width: 15 -> 20.5 (stored OOXML width). The stored "width" value
# derives from Excel's ColumnWidth (character units) via the sheet's max
# digit width; 59/3 characters round-trips to an OOXML width of 20.5 here.
$ws.Columns("A").ColumnWidth = 59 / 3

# Rows whose F value flips from 0 to 1
$rowsToFlag = @(2, 5, 7, 13, 15, 16, 23, 24, 26, 27, 33, 34, 36, 37, 38, 40, 42, 44, 46, 47, 48, 50, 51, 52, 53, 56, 57, 59, 60, 61, 62, 63, 64, 65, 66, 67, 68, 69, 70)

foreach ($r in $rowsToFlag) {
    $ws.Range("F$r").Value = 1
}

# Final selection/cursor position left by the author
$ws.Range("F71").Select()

# Best-effort: scroll the view so row 20 is the top visible row
# (matches the saved workbook's sheetView topLeftCell="A20").
$excel.ActiveWindow.ScrollRow = 20
